# Species-observation rows 2-6 on sheet "Artfynd" were re-sorted (the
# underlying records cycled positions: old row 2 -> row 3, old row 3 -> row 5,
# old row 4 -> row 6, old row 5 -> row 4, old row 6 -> row 2). Re-create that
# by writing each row's new content explicitly, clearing the few optional
# cells (J/L/M/AC/AF) that are empty/absent in their new position.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 111835838
$ws.Range("B2").Value = 89423
$ws.Range("E2").Value = 5432
$ws.Range("F2").Value = "Granticka"
$ws.Range("G2").Value = "Porodaedalea chrysoloma"
$ws.Range("H2").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("J2").Value = ""
$ws.Range("L2").ClearContents()
$ws.Range("M2").ClearContents()
$ws.Range("Q2").Value = 470914.6782613794
$ws.Range("R2").Value = 6810368.79402096
$ws.Range("S2").Value = 5
$ws.Range("AF2").Value = ""

# Row 3
$ws.Range("A3").Value = 111835718
$ws.Range("B3").Value = 56398
$ws.Range("E3").Value = 100109
$ws.Range("F3").Value = "Tretåig hackspett"
$ws.Range("G3").Value = "Picoides tridactylus"
$ws.Range("H3").Value = "(Linnaeus, 1758)"
$ws.Range("J3").ClearContents()
$ws.Range("L3").Value = ""
$ws.Range("M3").Value = "äldre spår"
$ws.Range("Q3").Value = 471101.0270993827
$ws.Range("R3").Value = 6810411.753755242
$ws.Range("S3").Value = 10
$ws.Range("AF3").ClearContents()

# Row 4
$ws.Range("A4").Value = 111835745
$ws.Range("B4").Value = 77515
$ws.Range("E4").Value = 6425
$ws.Range("F4").Value = "Garnlav"
$ws.Range("G4").Value = "Alectoria sarmentosa"
$ws.Range("H4").Value = "(Ach.) Ach."
$ws.Range("J4").Value = ""
$ws.Range("L4").ClearContents()
$ws.Range("M4").ClearContents()
$ws.Range("Q4").Value = 471152.5480076601
$ws.Range("R4").Value = 6810381.652036018
$ws.Range("AC4").ClearContents()
$ws.Range("AF4").Value = ""

# Row 5
$ws.Range("A5").Value = 111835758
$ws.Range("B5").Value = 77550
$ws.Range("E5").Value = 185
$ws.Range("F5").Value = "Violettgrå tagellav"
$ws.Range("G5").Value = "Bryoria nadvornikiana"
$ws.Range("H5").Value = "(Gyeln.) Brodo & D.Hawksw."
$ws.Range("Q5").Value = 471087.4311846643
$ws.Range("R5").Value = 6810390.807424263

# Row 6
$ws.Range("A6").Value = 111835826
$ws.Range("B6").Value = 56398
$ws.Range("E6").Value = 100109
$ws.Range("F6").Value = "Tretåig hackspett"
$ws.Range("G6").Value = "Picoides tridactylus"
$ws.Range("H6").Value = "(Linnaeus, 1758)"
$ws.Range("J6").ClearContents()
$ws.Range("L6").Value = ""
$ws.Range("M6").Value = "äldre spår"
$ws.Range("Q6").Value = 470915.776864712
$ws.Range("R6").Value = 6810385.536630718
$ws.Range("AC6").Value = "även hackspettbo, troligen av tret"
$ws.Range("AF6").ClearContents()
